# Auto-generated edit script: updates crypto price/volume data
# (and swaps the Avalanche/ShibaInu rows 15/16) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so Excel keeps values such as
# "1.00", "26.37" or "0.0000182" as literal text instead of coercing them
# into numbers (which would drop trailing zeros / use scientific notation).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.859.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.567.92'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.88'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.83'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.568.10'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +10.64%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.344'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.02'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.038.56'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000182'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.24%  '
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.37'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.773.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.563.55'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.69'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.20'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '366.13'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.09'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.80'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.30'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.23'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0922'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '517.60'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.79'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.41%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.25'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.07'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.91'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.76'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.95'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.326'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.06'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '152.61'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.63'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.524'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0260'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.62'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.42%  '
